# Apply the workbook edits described by the commit:
#   "add mythril bytecode results"
#
# This fills in previously-empty true/false positive/negative count cells
# (columns F:I) for the confuzzius tool rows, adds a "return success;"
# Property note to the mythril 0.7 block (L63:L66), tags M211 with the
# new "Unsafe Delegatecall" comment, and updates the sheet view
# (zoom/frozen pane/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- mythril 0.7 block: add "return success;" Property notes ---
$ws.Range("L63").Value = "return success;"
$ws.Range("L64").Value = "return success;"
$ws.Range("L65").Value = "return success;"
$ws.Range("L66").Value = "return success;"

# --- confuzzius tool blocks (versions 0.4, 0.5, 0.6, 0.7, 0.8): ---
# fill in the previously empty true positive / false positive /
# false negative / true negative counts for the "all", "delegatecall"
# (checked/unchecked) rows of each block.
$blockStarts = @(204, 224, 244, 264, 284)
foreach ($base in $blockStarts) {
    $rAll        = $base
    $rDcChecked  = $base + 2
    $rAllUnchk   = $base + 5
    $rDcUnchk    = $base + 7

    $ws.Range("F$rAll").Value = 0
    $ws.Range("G$rAll").Value = 3
    $ws.Range("H$rAll").Value = 0
    $ws.Range("I$rAll").Value = 32

    $ws.Range("F$rDcChecked").Value = 0
    $ws.Range("G$rDcChecked").Value = 0
    $ws.Range("H$rDcChecked").Value = 0
    $ws.Range("I$rDcChecked").Value = 10

    $ws.Range("F$rAllUnchk").Value = 12
    $ws.Range("G$rAllUnchk").Value = 0
    $ws.Range("H$rAllUnchk").Value = 16
    $ws.Range("I$rAllUnchk").Value = 0

    $ws.Range("F$rDcUnchk").Value = 0
    $ws.Range("G$rDcUnchk").Value = 0
    $ws.Range("H$rDcUnchk").Value = 8
    $ws.Range("I$rDcUnchk").Value = 0
}

# Comment on the delegatecall-unchecked row of the first confuzzius block
$ws.Range("M211").Value = "Unsafe Delegatecall"

# --- Sheet view: zoom + frozen pane scroll position + selection ---
$ws.Activate()
$ws.Range("A61").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("J73").Select() | Out-Null
$excel.ActiveWindow.Zoom = 77
